$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data
$ws.Range('D2').Value = "'67.316.49"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.78%  '
$ws.Range('D3').Value = "'3.251.57"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.55%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'577.76"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.39%  '
$ws.Range('D6').Value = "'179.12"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.56%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -0.58%  '
$ws.Range('D9').Value = "'3.249.63"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.61%  '
$ws.Range('E10').Value = '  +3.69%  '
$ws.Range('E11').Value = '  +1.76%  '
$ws.Range('D12').Value = "'0.414"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.44%  '
$ws.Range('D13').Value = "'3.813.22"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.50%  '
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('D15').Value = "'28.12"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.80%  '
$ws.Range('D16').Value = "'67.273.22"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.76%  '
$ws.Range('D17').Value = "'0.0000168"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.93%  '
$ws.Range('D18').Value = "'3.254.73"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.51%  '
$ws.Range('D19').Value = "'5.87"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.32%  '
$ws.Range('D20').Value = "'13.40"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.23%  '
$ws.Range('D21').Value = "'376.66"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.79%  '
$ws.Range('D22').Value = "'7.63"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.15%  '
$ws.Range('D23').Value = "'0.999"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').Value = "'71.30"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.48%  '
$ws.Range('D25').Value = "'0.512"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.73%  '
$ws.Range('D26').Value = "'3.395.79"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.51%  '
$ws.Range('D27').Value = "'0.0000118"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.25%  '
$ws.Range('D28').Value = "'9.99"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.56%  '
$ws.Range('E29').Value = '  +1.68%  '
$ws.Range('D30').Value = "'0.998"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('E31').Value = '  +4.19%  '
$ws.Range('D32').Value = "'5.64"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.41%  '
$ws.Range('D33').Value = "'22.61"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.66%  '
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('E35').Value = '  +5.55%  '
$ws.Range('D36').Value = "'6.85"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.24%  '
$ws.Range('D37').Value = "'163.48"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.08%  '
$ws.Range('E38').Value = '  +3.85%  '
$ws.Range('D39').Value = "'0.859"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.15%  '
$ws.Range('D40').Value = "'1.86"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.33%  '
$ws.Range('D41').Value = "'26.87"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.05%  '
$ws.Range('D42').Value = "'6.80"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +12.77%  '
$ws.Range('D43').Value = "'2.62"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.55%  '
$ws.Range('D44').Value = "'2.769.28"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.90%  '
$ws.Range('D45').Value = "'4.41"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.33%  '
$ws.Range('D46').Value = "'25.91"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.86%  '
$ws.Range('D47').Value = "'353.23"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +10.10%  '
$ws.Range('D48').Value = "'40.47"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.37%  '
$ws.Range('D49').Value = "'0.0675"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.14%  '
$ws.Range('D50').Value = "'0.0280"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.47%  '
$ws.Range('E51').Value = '  +1.50%  '
